# "latest pulled and increment 2 updated"
#
# Adds the Increment 2.2 (INC2.2) requirements-gathering tasks to the
# "Sprint 2" backlog sheet, fills in the previously-blank status cell for
# the last INC2.1 task, and leaves the workbook focused on the newly
# edited sheet/cell - matching the author's final view state.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sprint 2")

# --- finish off the INC2.1 row (E34 was blank, now marked complete) ---
$ws2.Range("E34").Value = "COM"

# --- new INC2.2 task rows (35-42) ---
$rows = @(
    @{ Row = 35; Assignee = "Cristian"; Description = "Game Controller Requirements" },
    @{ Row = 36; Assignee = "Haris";    Description = "Story Line Requirements" },
    @{ Row = 37; Assignee = "Bhuwan";   Description = "Audio and Visual Requirements" },
    @{ Row = 38; Assignee = "Diptin";   Description = "AI requirement" },
    @{ Row = 39; Assignee = "Cristian"; Description = "5) Software process and UML Diagrams" },
    @{ Row = 40; Assignee = "Bhuwan";   Description = "7) Delovery and Schedule" },
    @{ Row = 41; Assignee = "Haris";    Description = "6)Assumptions and constraints" },
    @{ Row = 42; Assignee = "Diptin";   Description = "Finalize and complete every requirements and combine" }
)

# carry the existing look-and-feel of the table down into the new rows
# (B34 = shaded/centered Task ID cell, E19 = shaded Status cell)
$ws2.Range("B34").Copy()
$ws2.Range("B35:B42").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("E19").Copy()
$ws2.Range("E35:E42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($r in $rows) {
    $row = $r.Row

    $ws2.Range("A" + $row).Formula = "=A" + ($row - 1) + "+1"
    $ws2.Range("B" + $row).Value = "INC2.2"
    $ws2.Range("C" + $row).Value = $r.Assignee
    $ws2.Range("D" + $row).Value = $r.Description
    $ws2.Range("E" + $row).Value = "x"
}

# --- workbook / view state: Sprint 2 becomes the active tab, cursor on G34 ---
$ws2.Activate()
$ws2.Range("G34").Select()
